# new card event balance
# Insert a new "QuestDungeonRate" table column between the existing
# "QuestDungeon" (L) and "BgImage" (M) columns, shifting BgImage to N,
# and populate the new column's header + per-row rate values. Also trim
# the "|bookancient;1" suffix that used to be tacked onto the
# QuestDungeon cell values, since that data now lives in its own column.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Insert a blank column at M (QuestDungeon stays at L, old M/BgImage
#    data + column formatting shifts right to N).
$ws.Range("M1").EntireColumn.Insert()

# 2) Grow the worksheet table to cover the new column.
$lo = $ws.ListObjects.Item(1)
$lo.Resize($ws.Range("A3:N7"))

# 3) Row 1 (merged/banner header row): mirror column L's banner cell into
#    the new column M (matches the source column being carried along by
#    the insert), keep N as the old BgImage banner cell (already shifted).
$ws.Range("M1").Value = $ws.Range("L1").Text

# 4) Row 3 (actual table header labels).
$ws.Range("M3").Value = "QuestDungeonRate"
$ws.Range("N3").Value = "BgImage"

# 5) Data rows: strip the trailing "|bookancient;1" from QuestDungeon (L)
#    and move the drop-rate information into the new QuestDungeonRate (M)
#    column.
$ws.Range("L4").Value = "trees;4|manflower;2|river;2|cliff;2|losttree;2|oldtree;1"
$ws.Range("M4").Value = "bookancient;1|zookeeper;1"

$ws.Range("L5").Value = "trees;4"
$ws.Range("M5").Value = ""

$ws.Range("L6").Value = "trees;2|sandland;2|potteryroom;2|honeyhome;2|snare;1|basement;1|woodhouse2;1|booty;1|trapspear;2|trapdrop;1|potteryman;1|stonedoor2;1|crystalball;2"
$ws.Range("M6").Value = "bookancient;1"

$ws.Range("L7").Value = "trees;4"
$ws.Range("M7").Value = "bookdead;2"
